$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("C2").Value = 4.92241824879848
$ws.Range("D2").Value = 3.721592044944652
$ws.Range("E2").Value = 16.57821045982143
$ws.Range("F2").Value = 19.9246460903633
$ws.Range("G2").Value = 21.66912798817862
$ws.Range("H2").Value = 12.18510546842633
$ws.Range("K2").Value = 12.35959968634289
$ws.Range("N2").Value = 15.89550549347478
$ws.Range("O2").Value = 17.63542988083426
$ws.Range("C3").Value = 4.749659201147663
$ws.Range("D3").Value = 3.667893492910732
$ws.Range("E3").Value = 15.63267006056087
$ws.Range("F3").Value = 19.86376230223798
$ws.Range("G3").Value = 21.50257220312157
$ws.Range("H3").Value = 12.22228320492159
$ws.Range("K3").Value = 11.72597441288274
$ws.Range("N3").Value = 15.91620507709672
$ws.Range("O3").Value = 17.65944936333058
$ws.Range("C4").Value = 4.64179640417249
$ws.Range("D4").Value = 3.634126685222137
$ws.Range("E4").Value = 15.02678396108317
$ws.Range("F4").Value = 19.8336983412799
$ws.Range("G4").Value = 21.41053999929148
$ws.Range("H4").Value = 12.24785591972771
$ws.Range("K4").Value = 11.34561334909411
$ws.Range("N4").Value = 15.93064810771888
$ws.Range("O4").Value = 17.67989248757084
$ws.Range("C5").Value = 4.597470968712622
$ws.Range("D5").Value = 3.620175410842088
$ws.Range("E5").Value = 14.77378080391067
$ws.Range("F5").Value = 19.8232945068944
$ws.Range("G5").Value = 21.37564753659426
$ws.Range("H5").Value = 12.2589653029393
$ws.Range("K5").Value = 11.18692944380788
$ws.Range("N5").Value = 15.93697050032259
$ws.Range("O5").Value = 17.68964915694793
$ws.Range("C6").Value = 4.590090829282954
$ws.Range("D6").Value = 3.61784754491931
$ws.Range("E6").Value = 14.73141028150578
$ws.Range("F6").Value = 19.82167869363531
$ws.Range("G6").Value = 21.37001237239117
$ws.Range("H6").Value = 12.26085151700205
$ws.Range("K6").Value = 11.16036347901897
$ws.Range("N6").Value = 15.93804673181066
$ws.Range("O6").Value = 17.69135517202771
$ws.Range("C7").Value = 4.641200006461131
$ws.Range("D7").Value = 3.633939294764685
$ws.Range("E7").Value = 15.02339617695049
$ws.Range("F7").Value = 19.83355054394364
$ws.Range("G7").Value = 21.4100588086384
$ws.Range("H7").Value = 12.24800296069571
$ws.Range("K7").Value = 11.34348793718519
$ws.Range("N7").Value = 15.93073160425642
$ws.Range("O7").Value = 17.6800183046552
$ws.Range("C8").Value = 4.863267485941556
$ws.Range("D8").Value = 3.703247267862178
$ws.Range("E8").Value = 16.25758812642456
$ws.Range("F8").Value = 19.90213938944155
$ws.Range("G8").Value = 21.60960042532906
$ws.Range("H8").Value = 12.1973533147878
$ws.Range("K8").Value = 12.1299998548495
$ws.Range("N8").Value = 15.90228346204964
$ws.Range("O8").Value = 17.64252633488655
$ws.Range("C9").Value = 5.281413424422786
$ws.Range("D9").Value = 3.832464971792832
$ws.Range("E9").Value = 18.58006853416102
$ws.Range("F9").Value = 20.09427424376053
$ws.Range("G9").Value = 22.08010339957468
$ws.Range("H9").Value = 12.11990319785334
$ws.Range("K9").Value = 13.69613409843147
$ws.Range("N9").Value = 15.86021306791182
$ws.Range("O9").Value = 17.61443053411613
$ws.Range("C10").Value = 5.574328220641153
$ws.Range("D10").Value = 3.922832740311731
$ws.Range("E10").Value = 20.23196380190635
$ws.Range("F10").Value = 20.26979644459703
$ws.Range("G10").Value = 22.47096386762556
$ws.Range("H10").Value = 12.07645443294198
$ws.Range("K10").Value = 14.73094982919387
$ws.Range("N10").Value = 15.8376160406507
$ws.Range("O10").Value = 17.62175201076088
$ws.Range("C11").Value = 5.703819808638455
$ws.Range("D11").Value = 3.962855152797303
$ws.Range("E11").Value = 20.94100851959499
$ws.Range("F11").Value = 20.35688717753232
$ws.Range("G11").Value = 22.65781841546636
$ws.Range("H11").Value = 12.05963366089517
$ws.Range("K11").Value = 15.17620792165485
$ws.Range("N11").Value = 15.82913013426447
$ws.Range("O11").Value = 17.63119357179574
$ws.Range("C12").Value = 5.75226807154351
$ws.Range("D12").Value = 3.977846788960541
$ws.Range("E12").Value = 21.20344337234398
$ws.Range("F12").Value = 20.39088565523874
$ws.Range("G12").Value = 22.72980636120647
$ws.Range("H12").Value = 12.05368926184906
$ws.Range("K12").Value = 15.34113124409634
$ws.Range("N12").Value = 15.82617371066759
$ws.Range("O12").Value = 17.63564916535606
$ws.Range("C13").Value = 5.741860686793722
$ws.Range("D13").Value = 3.97462548461605
$ws.Range("E13").Value = 21.14719254057301
$ws.Range("F13").Value = 20.38351853535839
$ws.Range("G13").Value = 22.71424892413851
$ws.Range("H13").Value = 12.05495055211255
$ws.Range("K13").Value = 15.30577627677552
$ws.Range("N13").Value = 15.82679901376327
$ws.Range("O13").Value = 17.63465040928804
$ws.Range("C14").Value = 5.707817698672183
$ws.Range("D14").Value = 3.964091853304557
$ws.Range("E14").Value = 20.96272055263594
$ws.Range("F14").Value = 20.35966396646909
$ws.Range("G14").Value = 22.66371667172097
$ws.Range("H14").Value = 12.05913607834218
$ws.Range("K14").Value = 15.18985024084266
$ws.Range("N14").Value = 15.8288817620525
$ws.Range("O14").Value = 17.63154248888423
$ws.Range("C15").Value = 5.686887586644592
$ws.Range("D15").Value = 3.957618120629455
$ws.Range("E15").Value = 20.84893728611592
$ws.Range("F15").Value = 20.34518435930422
$ws.Range("G15").Value = 22.63292224522514
$ws.Range("H15").Value = 12.06175526962762
$ws.Range("K15").Value = 15.11836161036336
$ws.Range("N15").Value = 15.83019094772652
$ws.Range("O15").Value = 17.62975346439383
$ws.Range("C16").Value = 5.565785881767384
$ws.Range("D16").Value = 3.920194648691373
$ws.Range("E16").Value = 20.18477552880464
$ws.Range("F16").Value = 20.26424875836209
$ws.Range("G16").Value = 22.45892858435091
$ws.Range("H16").Value = 12.07761311941208
$ws.Range("K16").Value = 14.70133634773959
$ws.Range("N16").Value = 15.83820662884686
$ws.Range("O16").Value = 17.62125810464613
$ws.Range("C17").Value = 5.490496310756023
$ws.Range("D17").Value = 3.896952673398997
$ws.Range("E17").Value = 19.76649539098715
$ws.Range("F17").Value = 20.21643756270825
$ws.Range("G17").Value = 22.35445749088073
$ws.Range("H17").Value = 12.08809686815229
$ws.Range("K17").Value = 14.43896004807847
$ws.Range("N17").Value = 15.84358273025155
$ws.Range("O17").Value = 17.61761305122429
$ws.Range("C18").Value = 5.446840572663826
$ws.Range("D18").Value = 3.883482692799758
$ws.Range("E18").Value = 19.52191782575629
$ws.Range("F18").Value = 20.18962125427619
$ws.Range("G18").Value = 22.29522482664659
$ws.Range("H18").Value = 12.09440386762291
$ws.Range("K18").Value = 14.28565020288317
$ws.Range("N18").Value = 15.84684379899251
$ws.Range("O18").Value = 17.61609167943261
$ws.Range("C19").Value = 5.432000707442147
$ws.Range("D19").Value = 3.878904745714722
$ws.Range("E19").Value = 19.43842108141276
$ws.Range("F19").Value = 20.18065974508947
$ws.Range("G19").Value = 22.27531885693344
$ws.Range("H19").Value = 12.09658682963986
$ws.Range("K19").Value = 14.23333104506649
$ws.Range("N19").Value = 15.84797697314458
$ws.Range("O19").Value = 17.61567529159162
$ws.Range("C20").Value = 5.498547718435722
$ws.Range("D20").Value = 3.899437415904225
$ws.Range("E20").Value = 19.81143513448486
$ws.Range("F20").Value = 20.22145657088921
$ws.Range("G20").Value = 22.36549049386462
$ws.Range("H20").Value = 12.08695217024751
$ws.Range("K20").Value = 14.46713880807355
$ws.Range("N20").Value = 15.84299296342262
$ws.Range("O20").Value = 17.61794153150099
$ws.Range("C21").Value = 5.717833239024674
$ws.Range("D21").Value = 3.967190346591713
$ws.Range("E21").Value = 21.01706878256086
$ws.Range("F21").Value = 20.36664316554223
$ws.Range("G21").Value = 22.67852641583445
$ws.Range("H21").Value = 12.05789513189421
$ws.Range("K21").Value = 15.22400071083317
$ws.Range("N21").Value = 15.82826304072356
$ws.Range("O21").Value = 17.63243146273351
$ws.Range("C22").Value = 5.857707254008134
$ws.Range("D22").Value = 4.010511657539388
$ws.Range("E22").Value = 21.76968408122538
$ws.Range("F22").Value = 20.46745902074983
$ws.Range("G22").Value = 22.89024879663818
$ws.Range("H22").Value = 12.04138422709254
$ws.Range("K22").Value = 15.69716181167201
$ws.Range("N22").Value = 15.82013376817003
$ws.Range("O22").Value = 17.64703240700331
$ws.Range("C23").Value = 5.783382922893529
$ws.Range("D23").Value = 3.987480486394014
$ws.Range("E23").Value = 21.37122132907077
$ws.Range("F23").Value = 20.413117369447
$ws.Range("G23").Value = 22.77662001873079
$ws.Range("H23").Value = 12.04996893608937
$ws.Range("K23").Value = 15.44659881462118
$ws.Range("N23").Value = 15.82433579175919
$ws.Range("O23").Value = 17.63876988347226
$ws.Range("C24").Value = 5.494908828247616
$ws.Range("D24").Value = 3.898314399162723
$ws.Range("E24").Value = 19.79113065880053
$ws.Range("F24").Value = 20.21918538544858
$ws.Range("G24").Value = 22.36049988496996
$ws.Range("H24").Value = 12.08746881684811
$ws.Range("K24").Value = 14.45440687078082
$ws.Range("N24").Value = 15.84325906651093
$ws.Range("O24").Value = 17.61779123696106
$ws.Range("C25").Value = 5.170560129404429
$ws.Range("D25").Value = 3.79827651602311
$ws.Range("E25").Value = 17.93391073534845
$ws.Range("F25").Value = 20.03619601637149
$ws.Range("G25").Value = 21.94463942545535
$ws.Range("H25").Value = 12.13850114885837
$ws.Range("K25").Value = 13.29264100259348
$ws.Range("N25").Value = 15.87013086466877
$ws.Range("O25").Value = 13.53211574969021
